$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.150.00"
$ws.Range("E2").Value = "  +0.05%  "
$ws.Range("D3").Value = "3.794.11"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'601.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.42%  "
$ws.Range("D6").Value = "'163.76"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.30%  "
$ws.Range("D7").Value = "3.792.31"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  -0.54%  "
$ws.Range("D10").Value = "'0.169"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.61%  "
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("E12").Value = "  -0.62%  "
$ws.Range("D13").Value = "'37.27"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.06%  "
$ws.Range("E14").Value = "  -2.26%  "
$ws.Range("D15").Value = "4.430.29"
$ws.Range("E15").Value = "  +1.12%  "
$ws.Range("D16").Value = "3.789.14"
$ws.Range("E16").Value = "  +0.88%  "
$ws.Range("D17").Value = "69.225.97"
$ws.Range("E17").Value = "  +0.20%  "
$ws.Range("E18").Value = "  +1.33%  "
$ws.Range("E19").Value = "  -0.26%  "
$ws.Range("D20").Value = "'17.27"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").Value = "'11.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.23%  "
$ws.Range("D22").Value = "'488.02"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.36%  "
$ws.Range("D23").Value = "'0.721"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("E24").Value = "  -2.75%  "
$ws.Range("D25").Value = "'84.60"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("E26").Value = "  -3.62%  "
$ws.Range("D27").Value = "'12.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.21%  "
$ws.Range("E28").Value = "  -2.75%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("E32").Value = "  -5.41%  "
$ws.Range("D33").Value = "3.945.55"
$ws.Range("E33").Value = "  +1.18%  "
$ws.Range("D34").Value = "'31.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.69%  "
$ws.Range("D35").Value = "3.742.01"
$ws.Range("E35").Value = "  +1.48%  "
$ws.Range("E36").Value = "  -2.15%  "
$ws.Range("E37").Value = "  +5.59%  "
$ws.Range("D38").Value = "'1.01"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  +0.05%  "
$ws.Range("E41").Value = "  -0.66%  "
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("D43").Value = "'48.60"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.04%  "
$ws.Range("E44").Value = "  +0.44%  "
$ws.Range("D45").Value = "'421.35"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.89%  "
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "'8.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.28%  "
$ws.Range("D48").Value = "2.832.56"
$ws.Range("E48").Value = "  +0.81%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'141.16"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("B50").Value = "Arweave"
$ws.Range("C50").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D50").Value = "'39.74"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").Value = "'0.0350"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.60%  "
